{"js": "// Replace the multiplication problems in the practice-sheet table with\n// newly generated ones (same \"AAA\u00d7B=\" shape), one-for-one, keeping all\n// run/paragraph formatting (font, size, etc.) intact.\nconst replacements = [\n  [\"871\u00d79=\", \"507\u00d73=\"],\n  [\"201\u00d78=\", \"424\u00d77=\"],\n  [\"638\u00d79=\", \"953\u00d78=\"],\n  [\"593\u00d79=\", \"892\u00d77=\"],\n  [\"582\u00d73=\", \"846\u00d78=\"],\n  [\"126\u00d74=\", \"493\u00d78=\"],\n  [\"323\u00d78=\", \"976\u00d74=\"],\n  [\"413\u00d76=\", \"994\u00d76=\"],\n  [\"468\u00d73=\", \"579\u00d77=\"],\n  [\"478\u00d73=\", \"523\u00d73=\"],\n  [\"642\u00d78=\", \"288\u00d75=\"],\n  [\"937\u00d74=\", \"674\u00d76=\"],\n  [\"696\u00d74=\", \"461\u00d73=\"],\n  [\"385\u00d75=\", \"767\u00d79=\"],\n  [\"494\u00d78=\", \"505\u00d75=\"],\n  [\"265\u00d75=\", \"683\u00d72=\"],\n  [\"316\u00d79=\", \"508\u00d75=\"],\n  [\"329\u00d76=\", \"882\u00d77=\"],\n  [\"902\u00d78=\", \"999\u00d77=\"],\n  [\"618\u00d77=\", \"431\u00d73=\"],\n  [\"515\u00d79=\", \"540\u00d75=\"],\n  [\"471\u00d77=\", \"858\u00d74=\"],\n  [\"777\u00d72=\", \"490\u00d73=\"],\n  [\"360\u00d75=\", \"177\u00d78=\"],\n  [\"436\u00d73=\", \"942\u00d73=\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Replace the multiplication problems in the practice-sheet table with\n# newly generated ones (same \"AAA\u00d7B=\" shape), one-for-one, keeping all\n# run/paragraph formatting (font, size, etc.) intact.\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"871\u00d79=\", \"507\u00d73=\"),\n    @(\"201\u00d78=\", \"424\u00d77=\"),\n    @(\"638\u00d79=\", \"953\u00d78=\"),\n    @(\"593\u00d79=\", \"892\u00d77=\"),\n    @(\"582\u00d73=\", \"846\u00d78=\"),\n    @(\"126\u00d74=\", \"493\u00d78=\"),\n    @(\"323\u00d78=\", \"976\u00d74=\"),\n    @(\"413\u00d76=\", \"994\u00d76=\"),\n    @(\"468\u00d73=\", \"579\u00d77=\"),\n    @(\"478\u00d73=\", \"523\u00d73=\"),\n    @(\"642\u00d78=\", \"288\u00d75=\"),\n    @(\"937\u00d74=\", \"674\u00d76=\"),\n    @(\"696\u00d74=\", \"461\u00d73=\"),\n    @(\"385\u00d75=\", \"767\u00d79=\"),\n    @(\"494\u00d78=\", \"505\u00d75=\"),\n    @(\"265\u00d75=\", \"683\u00d72=\"),\n    @(\"316\u00d79=\", \"508\u00d75=\"),\n    @(\"329\u00d76=\", \"882\u00d77=\"),\n    @(\"902\u00d78=\", \"999\u00d77=\"),\n    @(\"618\u00d77=\", \"431\u00d73=\"),\n    @(\"515\u00d79=\", \"540\u00d75=\"),\n    @(\"471\u00d77=\", \"858\u00d74=\"),\n    @(\"777\u00d72=\", \"490\u00d73=\"),\n    @(\"360\u00d75=\", \"177\u00d78=\"),\n    @(\"436\u00d73=\", \"942\u00d73=\")\n)\n\nforeach ($pair in $replacements) {\n    $oldText = $pair[0]\n    $newText = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $oldText\n    $find.Replacement.Text = $newText\n    $find.Forward = $true\n    $find.Wrap = 1\n    $find.Execute($find.Text, $false, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)\n}\n"}
